$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddProduct")

# Update the image filename used for the product in row 4 (L4)
$ws.Range("L4").Value = "gio-qua-2.jpg"

# Update the active selection on the sheet to L6
$ws.Range("L6").Select()
